# Auto-generated edit script applying the Universalis market-data refresh diff
# to the Kujata_Profits workbook (per-sheet Leve profit columns H..N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6958793.5
$ws.Range("I43").Value = 18965.166
$ws.Range("J43").Value = 27778278
$ws.Range("K43").Value = 18965.166
$ws.Range("L43").Value = 27778278
$ws.Range("M43").Value = -18896.166
$ws.Range("N43").Value = -27778416
$ws.Range("H112").Value = 2312.818
$ws.Range("J112").Value = 2660.1482
$ws.Range("L112").Value = 7980.444600000001
$ws.Range("N112").Value = -10196.4446
$ws.Range("H132").Value = 7944383.5
$ws.Range("I132").Value = 11910778
$ws.Range("K132").Value = 35732334
$ws.Range("M132").Value = -35729804
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3362.9092
$ws.Range("I32").Value = 3007.1746
$ws.Range("K32").Value = 3007.1746
$ws.Range("M32").Value = -2720.1746
$ws.Range("H63").Value = 1182
$ws.Range("I63").Value = 1002.5
$ws.Range("J63").Value = 1900
$ws.Range("K63").Value = 1002.5
$ws.Range("L63").Value = 1900
$ws.Range("M63").Value = -316.5
$ws.Range("N63").Value = -3272
$ws.Range("H66").Value = 1182
$ws.Range("I66").Value = 1002.5
$ws.Range("J66").Value = 1900
$ws.Range("K66").Value = 5012.5
$ws.Range("L66").Value = 9500
$ws.Range("M66").Value = -1580.5
$ws.Range("N66").Value = -16364
$ws.Range("H110").Value = 1250.7391
$ws.Range("I110").Value = 1076.2858
$ws.Range("J110").Value = 1522.1111
$ws.Range("K110").Value = 1076.2858
$ws.Range("L110").Value = 1522.1111
$ws.Range("M110").Value = 968.7141999999999
$ws.Range("N110").Value = -5612.1111
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1309.8572
$ws.Range("I80").Value = 892.3333
$ws.Range("K80").Value = 892.3333
$ws.Range("M80").Value = 105.6667
$ws.Range("H83").Value = 1309.8572
$ws.Range("I83").Value = 892.3333
$ws.Range("K83").Value = 4461.6665
$ws.Range("M83").Value = 530.3334999999997
$ws.Range("H134").Value = 4929.6763
$ws.Range("I134").Value = 1138.1482
$ws.Range("J134").Value = 19554.143
$ws.Range("K134").Value = 3414.4446
$ws.Range("L134").Value = 58662.429
$ws.Range("M134").Value = -879.4446000000003
$ws.Range("N134").Value = -63732.429
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 5587.25
$ws.Range("J43").Value = 5587.25
$ws.Range("L43").Value = 5587.25
$ws.Range("N43").Value = -5955.25
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H101").Value = 5587.25
$ws.Range("J101").Value = 5587.25
$ws.Range("L101").Value = 5587.25
$ws.Range("N101").Value = -12077.25
$ws.Range("H134").Value = 899.4286
$ws.Range("I134").Value = 918.6061
$ws.Range("K134").Value = 2755.8183
$ws.Range("M134").Value = -220.8182999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1239.9722
$ws.Range("J5").Value = 792.7778
$ws.Range("L5").Value = 2378.3334
$ws.Range("N5").Value = -2602.3334
$ws.Range("H12").Value = 64.82143000000001
$ws.Range("I12").Value = 76.25
$ws.Range("J12").Value = 60.25
$ws.Range("K12").Value = 228.75
$ws.Range("L12").Value = 180.75
$ws.Range("M12").Value = -55.75
$ws.Range("N12").Value = -526.75
$ws.Range("H94").Value = 4417.7144
$ws.Range("I94").Value = 5024
$ws.Range("J94").Value = 4316.6665
$ws.Range("K94").Value = 15072
$ws.Range("L94").Value = 12949.9995
$ws.Range("M94").Value = -14396
$ws.Range("N94").Value = -14301.9995
$ws.Range("H131").Value = 23810742
$ws.Range("J131").Value = 1722.7142
$ws.Range("L131").Value = 5168.142599999999
$ws.Range("N131").Value = -15248.1426
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 31500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -28970
$ws.Range("N132").Value = $null
$ws.Range("H135").Value = 1239.9722
$ws.Range("J135").Value = 792.7778
$ws.Range("L135").Value = 7135.000199999999
$ws.Range("N135").Value = -12205.0002
$ws.Range("H139").Value = 1903.92
$ws.Range("I139").Value = 1968.3158
$ws.Range("J139").Value = 1700
$ws.Range("K139").Value = 5904.9474
$ws.Range("L139").Value = 5100
$ws.Range("M139").Value = -764.9474
$ws.Range("N139").Value = -15380
$ws.Range("H140").Value = 22718.182
$ws.Range("I140").Value = 52424.453
$ws.Range("J140").Value = 2914
$ws.Range("K140").Value = 157273.359
$ws.Range("L140").Value = 8742
$ws.Range("M140").Value = -152093.359
$ws.Range("N140").Value = -19102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 30784.875
$ws.Range("J86").Value = 30784.875
$ws.Range("L86").Value = 30784.875
$ws.Range("N86").Value = -33156.875
$ws.Range("H89").Value = 30784.875
$ws.Range("J89").Value = 30784.875
$ws.Range("L89").Value = 92354.625
$ws.Range("N89").Value = -104210.625
$ws.Range("H102").Value = 7071.6816
$ws.Range("I102").Value = 5882.0713
$ws.Range("J102").Value = 9153.5
$ws.Range("K102").Value = 5882.0713
$ws.Range("L102").Value = 9153.5
$ws.Range("M102").Value = -4260.0713
$ws.Range("N102").Value = -12397.5
$ws.Range("H132").Value = 1964.2858
$ws.Range("I132").Value = 1408.091
$ws.Range("J132").Value = 4003.6667
$ws.Range("K132").Value = 4224.272999999999
$ws.Range("L132").Value = 12011.0001
$ws.Range("M132").Value = -1694.272999999999
$ws.Range("N132").Value = -17071.0001
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 359.58334
$ws.Range("I55").Value = 279.6154
$ws.Range("J55").Value = 454.0909
$ws.Range("K55").Value = 279.6154
$ws.Range("L55").Value = 454.0909
$ws.Range("M55").Value = -106.6154
$ws.Range("N55").Value = -800.0908999999999
$ws.Range("H132").Value = 21935.8
$ws.Range("I132").Value = 1491.6538
$ws.Range("K132").Value = 4474.9614
$ws.Range("M132").Value = -1944.9614
$ws.Range("H136").Value = 1387.56
$ws.Range("I136").Value = 1178.421
$ws.Range("J136").Value = 2049.8333
$ws.Range("K136").Value = 3535.263
$ws.Range("L136").Value = 6149.499899999999
$ws.Range("M136").Value = -985.2629999999999
$ws.Range("N136").Value = -11249.4999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 486.25
$ws.Range("I113").Value = 359.44446
$ws.Range("K113").Value = 1078.33338
$ws.Range("M113").Value = 1091.66662
$ws.Range("H132").Value = 1882.0927
$ws.Range("I132").Value = 1917.6086
$ws.Range("J132").Value = 1677.875
$ws.Range("K132").Value = 5752.825800000001
$ws.Range("L132").Value = 5033.625
$ws.Range("M132").Value = -3222.825800000001
$ws.Range("N132").Value = -10093.625
$ws.Range("H136").Value = 544.2381
$ws.Range("I136").Value = 267.23077
$ws.Range("J136").Value = 994.375
$ws.Range("K136").Value = 801.69231
$ws.Range("L136").Value = 2983.125
$ws.Range("M136").Value = 1748.30769
$ws.Range("N136").Value = -8083.125
